# Add a new "september-2025" worksheet at the end of the workbook,
# mirroring the existing monthly sheets (e.g. "august-2025"): a single
# cell A1 containing the month's Corporation Tax summary line.

$wb = $excel.ActiveWorkbook

$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)

# Insert the new sheet right after the current last sheet so it lands
# at the end of the tab order (like the diff's appended <sheet> entry).
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "september-2025"

$ws.Range("A1").Value = "Corporation Tax                               18,246             17,804                   443                2.5%"
